$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Mark "Nylon Bearing" (row 3) and "Threaded rod" (row 4) as Acquired
$ws.Range("D3").Value = "x"
$ws.Range("D4").Value = "x"

# Add descriptions for "9396K115" (row 11) and "4452K141" (row 20)
$ws.Range("B11").Value = "O-ring"
$ws.Range("B20").Value = "M10 x1"

# Widen column B to fit the new descriptions
# (ColumnWidth is offset from the stored OOXML width by 5/6, so back that out
#  to land exactly on the target stored width of 20.5)
$ws.Columns.Item(2).ColumnWidth = 20.5 - 5/6

# Update the active selection to B21
$ws.Range("B21").Select()
